# "Changes done for Kaman new UI - header & footer"
#  - Testdata sheet: two new key/value rows (EleType1/EleType2 -> JSElement)
#  - Selection/active-cell bookkeeping on both worksheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC07_SearchCategory")
$ws2 = $wb.Worksheets.Item("Testdata")

# --- Testdata sheet: append new key/value rows -----------------------
# Write column A first (creates the EleType1 / EleType2 shared strings),
# then column B (creates the shared JSElement string) so the shared
# string table ends up in the same order as the target workbook.
$ws2.Range("A15").Value = "EleType1"
$ws2.Range("A16").Value = "EleType2"
$ws2.Range("B15").Value = "JSElement"
$ws2.Range("B16").Value = "JSElement"

# give the new rows the same thin-box border used by the rest of the table
$newRows = $ws2.Range("A15:B16")
$newRows.Borders.LineStyle = 1
$newRows.Borders.Weight = 2

# --- Selections shown by the saved workbook ---------------------------
$ws1.Activate()
$ws1.Range("A3:XFD7").Select()

$ws2.Activate()
$ws2.Range("A15:B16").Select()

# leave the first sheet as the active/visible tab, matching the workbook
$ws1.Activate()
